$d = $word.ActiveDocument

# --- 1) Insert a new "Meta description" paragraph right after the first
#        (Heading1) paragraph, mirroring the empty-run + bold-run +
#        plain-run structure used elsewhere in the document. ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"

$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
           "<w:r/>" + `
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" + `
           "<w:r><w:t>: Read our review of Bigger Bass Blizzard - Christmas Catch and play it for free. Enjoy stunning graphics, high maximum win, and a special bonus game.</w:t></w:r>" + `
           "</w:p>"
$newPara.Range.InsertXML($metaXml)

# --- 2) Near the end of the document, drop the duplicated bold
#        "Play Bigger Bass Blizzard - Christmas Catch Free | Review"
#        paragraph entirely, and replace the meta-description text of the
#        final (italic) paragraph with the new AI image prompt, keeping
#        its italic formatting. These are now the last two paragraphs. ---
$n = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($n - 1)
$dupTitlePara.Range.Delete()

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$lastPara.Range.Find.ClearFormatting()
$lastPara.Range.Find.Execute(
    "Read our review of Bigger Bass Blizzard - Christmas Catch and play it for free. Enjoy stunning graphics, high maximum win, and a special bonus game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create an image featuring a happy Maya warrior with glasses in a cartoon style. The warrior should have a festive look, with a Santa hat, a scarf, and a fishing rod in hand, ready to catch some big bass in the frozen lake. In the background, there should be snow-covered trees and white flakes falling, creating a perfect Christmas atmosphere. The image should have bright and colorful tones to make it eye-catching and appealing to the players. The goal is to showcase the fun and thrilling experience of the game, while also highlighting the festive season and the unique character of the Maya warrior.",
    2)

Write-Output "edit complete"
